# [SAU-1563] combine ices div and subdiv
#
# The worksheet has two adjacent header columns, "ICES division" (H) and
# "ICES subdivision" (I). This change combines them into a single column,
# "ICES area", in column H, and removes the old "ICES subdivision" column
# (I) entirely - shifting every column to its right one position to the
# left (J->I, K->J, ... AB->AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "ICES division" header to the combined "ICES area" label.
$ws.Range("H1").Value = "ICES area"

# Remove the now redundant "ICES subdivision" column; Excel automatically
# shifts everything to the right of it one column to the left.
$ws.Columns("I").Delete()
